$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two time-range labels in column B
$ws.Range("B8").Value = "07:45 - 07:49"
$ws.Range("B9").Value = "07:50 - 07:54"

# Update the active selection on the sheet
$ws.Range("A15:B15").Select()
$ws.Range("B15").Activate()
